$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = -21.11380000000001
$ws.Range("A10").Value = -20.46549999999997
$ws.Range("A12").Value = -22.41730000000004
$ws.Range("D13").Value = -7.940100000000005
$ws.Range("A18").Value = -22.24880000000003
